$d = $word.ActiveDocument

# Locate the paragraph that contains exactly "Introduction" (the section heading),
# then operate on the very next paragraph (the empty placeholder paragraph that
# currently carries bullet/numbering formatting and needs to receive the new text).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Introduction") {
        $target = $p.Next()
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the placeholder paragraph following 'Introduction'."
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="onvn"/><w:ind w:left="720"/></w:pPr><w:r><w:t>On average, each semester a student will have to study around 7-8 subjects, equivalent to 7-8 textbooks.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Among them are many textbooks used only once.</w:t></w:r><w:r><w:t xml:space="preserve"> “</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>H</w:t></w:r><w:r><w:t>and</w:t></w:r><w:r><w:t>B</w:t></w:r><w:r><w:t>ook</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>serves as a marketplace for students to</w:t></w:r><w:r><w:t xml:space="preserve"> sell</w:t></w:r><w:r><w:t xml:space="preserve"> old</w:t></w:r><w:r><w:t xml:space="preserve"> materials and learning tools </w:t></w:r><w:r><w:t>they no longer need to other students</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>This helps us make the most of available resources.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>The products can be old books, old documents, old calculators, etc. In addition, this website also welcomes other sellers like bookstores, stationery stores, printing shops.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)
